$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.683.72"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "1.828.05"
$ws.Range("E3").Value = "  +1.79%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.47"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4660"
$ws.Range("E7").Value = "  +3.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3604"
$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07136"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9043"
$ws.Range("E10").Value = "  +2.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07771"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.46"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("D13").Value = "1.834.08"
$ws.Range("E13").Value = "  +1.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.267"
$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.343"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.79"
$ws.Range("E16").Value = "  +3.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008571"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").Value = "26.724.03"
$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.920"
$ws.Range("E24").Value = "  -2.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.45"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.93"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("E27").Value = "  -2.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.96"
$ws.Range("E28").Value = "  +1.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.833"
$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08807"
$ws.Range("E30").Value = "  +1.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.149"
$ws.Range("E31").Value = "  +2.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.762"
$ws.Range("E32").Value = "  +0.61%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7348"
$ws.Range("E33").Value = "  +1.67%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.153"
$ws.Range("E34").Value = "  +4.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.441"
$ws.Range("E35").Value = "  -0.25%  "

$ws.Range("E36").Value = "  +0.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01927"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.929"
$ws.Range("E38").Value = "  +2.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05133"
$ws.Range("E39").Value = "  +0.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.878"
$ws.Range("E40").Value = "  -0.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5073"
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1501"
$ws.Range("E42").Value = "  -0.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.040"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("E44").Value = "  +1.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.007"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.982"
$ws.Range("E46").Value = "  +1.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.32"
$ws.Range("E47").Value = "  -2.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.566"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06060"
$ws.Range("E49").Value = "  +1.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.03"
$ws.Range("E50").Value = "  -0.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.90"
$ws.Range("E51").Value = "  -0.37%  "
